# NATMI LR-pairs output (Cxcl5-Cxcr2) regenerated against the updated TPM
# matrix. The "Sending cluster" label "MuSCs" is renamed to "Resolving-Mac"
# (which collapses it with the existing "Resolving-Mac" target-cluster
# label), and every NATMI-derived statistic column (specificity scores,
# expression weights, edge weights, etc.) is recomputed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "MuSCs" sending-cluster label to "Resolving-Mac" ----------
# (rows 4-5, column A) and keep the existing "Resolving-Mac" target-cluster
# label (rows 3-5, column D) as-is/explicit.
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("D5").Value = "Resolving-Mac"

# --- Row 2: FAPs -> Cxcl5/Cxcr2 -> ECs -------------------------------------
$ws.Range("I2").Value = 0.9961712500318616
$ws.Range("J2").Value = 0.9961712500318616
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 0.1847226071284444
$ws.Range("R2").Value = 1.662503464156
$ws.Range("S2").Value = 0.7404068059287783
$ws.Range("T2").Value = 0.7404068059287784

# --- Row 3: FAPs -> Cxcl5/Cxcr2 -> Resolving-Mac ---------------------------
$ws.Range("I3").Value = 0.9961712500318616
$ws.Range("J3").Value = 0.9961712500318616
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 0.06381015753388888
$ws.Range("R3").Value = 0.5742914178049999
$ws.Range("S3").Value = 0.2557644441030832
$ws.Range("T3").Value = 0.2557644441030832

# --- Row 4: Resolving-Mac (was MuSCs) -> Cxcl5/Cxcr2 -> ECs ----------------
$ws.Range("G4").Value = 0.01807533333333334
$ws.Range("H4").Value = 0.054226
$ws.Range("I4").Value = 0.003828749968138469
$ws.Range("J4").Value = 0.003828749968138468
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 0.0007099749928888889
$ws.Range("R4").Value = 0.006389774936000001
$ws.Range("S4").Value = 0.002845728116042946
$ws.Range("T4").Value = 0.002845728116042945

# --- Row 5: Resolving-Mac (was MuSCs) -> Cxcl5/Cxcr2 -> Resolving-Mac ------
$ws.Range("G5").Value = 0.01807533333333334
$ws.Range("H5").Value = 0.054226
$ws.Range("I5").Value = 0.003828749968138469
$ws.Range("J5").Value = 0.003828749968138468
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 0.0002452521477777778
$ws.Range("R5").Value = 0.00220726933
$ws.Range("S5").Value = 0.0009830218520955234
$ws.Range("T5").Value = 0.0009830218520955232
